$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-29 Tuesday" "2023-08-30 Wednesday"

Replace-Text "98×59=" "86×84="
Replace-Text "52×13=" "41×23="
Replace-Text "51×63=" "42×61="
Replace-Text "43×25=" "16×42="
Replace-Text "82×32=" "77×76="
Replace-Text "30×58=" "23×64="
Replace-Text "39×52=" "67×85="
Replace-Text "95×18=" "64×58="
Replace-Text "31×30=" "47×43="
Replace-Text "14×16=" "69×99="
Replace-Text "69×28=" "48×87="
Replace-Text "94×19=" "40×92="
Replace-Text "64×16=" "16×15="
Replace-Text "40×88=" "92×96="
Replace-Text "23×13=" "48×53="
Replace-Text "97×22=" "59×56="
Replace-Text "40×69=" "25×87="
Replace-Text "57×25=" "25×18="
Replace-Text "70×59=" "73×46="
Replace-Text "33×73=" "27×56="
Replace-Text "51×84=" "35×75="
Replace-Text "88×26=" "29×90="
Replace-Text "56×84=" "58×86="
Replace-Text "47×96=" "64×35="
Replace-Text "89×30=" "68×12="
